$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I13").Value = 0.3779298801345128
$ws.Range("J13").Value = 0.08732206993011483
$ws.Range("K13").Value = -0.3326748274291247
$ws.Range("L13").Value = 1.552525783804963
